# Update crypto price/volume data (and swap Aave/Frax row order) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.890.07'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.635.94'
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.38'
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5067'
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2573'
$ws.Range("E8").Value = '  +0.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06358'
$ws.Range("E9").Value = '  -0.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.72'
$ws.Range("E10").Value = '  +1.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07754'
$ws.Range("E11").Value = '  -0.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.289'
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.619.21'
$ws.Range("E13").Value = '  -1.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5449'
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0₅7736'
$ws.Range("E15").Value = '  -1.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.10'
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.922.53'
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.442'
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '195.68'
$ws.Range("E20").Value = '  -1.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.923'
$ws.Range("E21").Value = '  -0.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.131'
$ws.Range("E22").Value = '  +1.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.004'
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("E24").Value = '  +0.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.95'
$ws.Range("E25").Value = '  +1.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1245'
$ws.Range("E26").Value = '  +8.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.836'
$ws.Range("E27").Value = '  -0.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.63'
$ws.Range("E28").Value = '  -0.70%  '
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.04872'
$ws.Range("E30").Value = '  -3.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.240'
$ws.Range("E31").Value = '  -0.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.197'
$ws.Range("E32").Value = '  +0.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.550'
$ws.Range("E33").Value = '  +0.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.372'
$ws.Range("E34").Value = '  +0.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9124'
$ws.Range("E35").Value = '  +1.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.571'
$ws.Range("E36").Value = '  -1.01%  '
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.122.17'
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01564'
$ws.Range("E39").Value = '  +0.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.002'
$ws.Range("E40").Value = '  -0.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.591'
$ws.Range("E41").Value = '  -0.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8043'
$ws.Range("E42").Value = '  -1.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '98.55'
$ws.Range("E43").Value = '  -1.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₈123'
$ws.Range("E44").Value = '  -7.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.769.37'
$ws.Range("E45").Value = '  -0.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4483'
$ws.Range("E46").Value = '  -1.04%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.006'
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.01'
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05182'
$ws.Range("E49").Value = '  +2.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.531'
$ws.Range("E50").Value = '  +1.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.003'
$ws.Range("E51").Value = '  -0.33%  '
